$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 ("grandes regiões e unidades da federação") is an empty section-header
# row with no data. Delete it entirely; all rows below (7..37) shift up by one,
# and the now-unused shared string is dropped automatically on save.
$ws.Rows(6).Delete()
